$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet, positioned right after "总计" and
#    before "2021-Q4" (i.e. as the 2nd tab).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$oldQ4_2021 = $wb.Worksheets.Item(2)
$wb.Worksheets.Add($oldQ4_2021) | Out-Null
$wb.Worksheets.Item(2).Name = "2022-Q3"

# Re-fetch sheets by name (fresh references) after the rename.
$q3_2022 = $wb.Worksheets.Item("2022-Q3")
$q4_2021 = $wb.Worksheets.Item("2021-Q4")

# --- Header row (copy formatting from the sibling "2021-Q4" sheet so the
#     new sheet reuses the same shared style, then overwrite the text). ---
$q4_2021.Range("B1:H1").Copy()
$q3_2022.Range("B1:H1").PasteSpecial(-4122)

$q3_2022.Range("B1").Value = "基金代码"
$q3_2022.Range("C1").Value = "基金名称"
$q3_2022.Range("D1").Value = "基金规模"
$q3_2022.Range("E1").Value = "股票总仓位"
$q3_2022.Range("F1").Value = "仓位占比"
$q3_2022.Range("G1").Value = "持有市值(亿元)"
$q3_2022.Range("H1").Value = "仓位排名"

# --- Data row 2 ---
$q4_2021.Range("A2").Copy()
$q3_2022.Range("A2").PasteSpecial(-4122)
$q3_2022.Range("A2").Value = 0

# Fund code / name / figures are stored as plain text in this workbook
# (even though several look numeric), so force text via a temporary "@"
# number format, then clear the leftover formatting so the cell ends up
# with the default style (matching the rest of the sheet).
$q3_2022.Range("B2").NumberFormat = "@"
$q3_2022.Range("B2").Value = "002952"
$q3_2022.Range("B2").ClearFormats()

$q3_2022.Range("C2").Value = "建信多因子量化股票"

$q3_2022.Range("D2").NumberFormat = "@"
$q3_2022.Range("D2").Value = "0.09"
$q3_2022.Range("D2").ClearFormats()

$q3_2022.Range("E2").NumberFormat = "@"
$q3_2022.Range("E2").Value = "91.26"
$q3_2022.Range("E2").ClearFormats()

$q3_2022.Range("F2").NumberFormat = "@"
$q3_2022.Range("F2").Value = "3.01"
$q3_2022.Range("F2").ClearFormats()

$q3_2022.Range("G2").NumberFormat = "@"
$q3_2022.Range("G2").Value = "0.0027"
$q3_2022.Range("G2").ClearFormats()

$q3_2022.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: shift the existing 3 rows down by one
#    and insert the new "2022-Q3" row at the top of the data (row 2).
# ---------------------------------------------------------------------------
$ws = $total

# Extend the styled index column (A) down to the new last row (row 5) by
# copying the format already present on A4.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = "2021-Q4"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0.04

$ws.Range("B4").Value = "2021-Q2"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0.14

$ws.Range("B5").Value = "2020-Q4"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 4.67
